$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.620.58'
$ws.Range('E2').Value = '  +3.11%  '
$ws.Range('D3').Value = '2.446.08'
$ws.Range('E3').Value = '  +1.97%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.46'
$ws.Range('E5').Value = '  +2.69%  '
$ws.Range('E6').Value = '  +2.92%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +0.52%  '
$ws.Range('D9').Value = '2.445.20'
$ws.Range('E9').Value = '  +1.65%  '
$ws.Range('E10').Value = '  +1.95%  '
$ws.Range('E11').Value = '  +0.91%  '
$ws.Range('B12').Value = 'Toncoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.23'
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('B13').Value = 'Cardano'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.353'
$ws.Range('E13').Value = '  +3.06%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.44'
$ws.Range('E14').Value = '  +8.89%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000178'
$ws.Range('E15').Value = '  +5.39%  '
$ws.Range('B16').Value = 'BabyDogeCoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D16').Value = '0.0₅0128'
$ws.Range('E16').Value = '  +356.33%  '
$ws.Range('D17').Value = '2.889.48'
$ws.Range('E17').Value = '  +1.88%  '
$ws.Range('D18').Value = '62.540.18'
$ws.Range('E18').Value = '  +3.25%  '
$ws.Range('D19').Value = '2.448.14'
$ws.Range('E19').Value = '  +1.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.79'
$ws.Range('E20').Value = '  -2.64%  '
$ws.Range('E21').Value = '  +2.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '326.51'
$ws.Range('E22').Value = '  +0.80%  '
$ws.Range('E24').Value = '  +10.68%  '
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '65.43'
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '638.64'
$ws.Range('E27').Value = '  +13.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.16'
$ws.Range('E28').Value = '  +15.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.53'
$ws.Range('E29').Value = '  +5.88%  '
$ws.Range('D30').Value = '0.0₃0978'
$ws.Range('E30').Value = '  +4.47%  '
$ws.Range('D31').Value = '2.567.11'
$ws.Range('E31').Value = '  +2.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.17'
$ws.Range('E32').Value = '  +1.08%  '
$ws.Range('E33').Value = '  +6.34%  '
$ws.Range('E34').Value = '  +3.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.138'
$ws.Range('E35').Value = '  +5.69%  '
$ws.Range('E36').Value = '  +1.77%  '
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('E38').Value = '  +3.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '153.72'
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.45'
$ws.Range('E40').Value = '  +6.20%  '
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('E42').Value = '  +1.68%  '
$ws.Range('E43').Value = '  +7.84%  '
$ws.Range('E44').Value = '  +4.94%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.57'
$ws.Range('E45').Value = '  +2.12%  '
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('E47').Value = '  +28.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '144.32'
$ws.Range('E48').Value = '  +1.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.58'
$ws.Range('E49').Value = '  +1.31%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.50'
$ws.Range('E50').Value = '  +6.24%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.603'
$ws.Range('E51').Value = '  +2.19%  '
